$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "X axis"
$ws.Range("C8").Value = "Added start and end fields."
$ws.Range("D8").Value = "1. Enter value for Start field" + [char]10 + "2. Enter value for End field"
$ws.Range("E8").Value = "Histogram Chart will be displayed with specified start and end values."

$ws.Range("C7:E7").Copy()
$ws.Range("C8:E8").PasteSpecial(-4122)
$ws.Rows.Item(8).AutoFit()

$ws.Range("E9").Select()
